$d = $word.ActiveDocument

# ------------------------------------------------------------------
# 1) Insert two new body paragraphs right after the "Write Up" title
#    paragraph (paragraph 1). These become paragraphs 2 and 3.
# ------------------------------------------------------------------
$titlePara = $d.Paragraphs(1)
$titlePara.Range.InsertParagraphAfter()
$d.Paragraphs(2).Range.InsertParagraphAfter()

# --- paragraph 2: the "This week, we will be looking..." blurb -----
$introText = "This week, we will be looking at how to use mixins with media queries. If you have different devices, you will need your layout to display at different widths for each of these different formats. In this tutorial, we will be writing some mixin code to easily switch between the layouts of desktop and cell phone, and it even handles the unique sizing for the tablet as well. "

$p2 = $d.Paragraphs(2)
$p2.Style = "Normal"
$p2.Range.Text = $introText

# Split the run containing "mixins" into three runs (matching the
# word being wrapped in spell-check markers in the authored copy) by
# bookmarking the word and then removing the bookmark again -- this
# forces a run boundary at each edge of the word without leaving any
# bookmark behind.
$p2Start = $d.Paragraphs(2).Range.Start
$mixinsOffset = $introText.IndexOf("mixins")
$mixinsStart = $p2Start + $mixinsOffset
$mixinsEnd = $mixinsStart + "mixins".Length
$mixinsRange = $d.Range($mixinsStart, $mixinsEnd)
$d.Bookmarks.Add("zzz_mixins_split", $mixinsRange)
$d.Bookmarks("zzz_mixins_split").Delete()

# --- paragraph 3: the "So, if this is something..." sentence -------
$p3 = $d.Paragraphs(3)
$p3.Style = "Normal"
$p3.Range.Text = "So, if this is something that you would like to take a more in-depth look at, then please join us for our new article this week entitled:"

# ------------------------------------------------------------------
# 2) Insert the new Heading1 paragraph right after the first of the
#    pre-existing empty paragraphs (now paragraph 4).
# ------------------------------------------------------------------
$emptyPara = $d.Paragraphs(4)
$emptyPara.Range.InsertParagraphAfter()

$headingPara = $d.Paragraphs(5)
$headingPara.Style = "Heading1"
$headingPara.Range.Text = "9 Using a Mixin for Media Queries"

Write-Output "Final paragraph count: $($d.Paragraphs.Count)"
foreach ($p in $d.Paragraphs) {
    Write-Output ("[" + $p.Style.NameLocal + "] " + $p.Range.Text)
}
